$d = $word.ActiveDocument

# Remove " with a high degree of accuracy" before the period following
# "...detect moustaches within images of faces"
$d.Content.Find.Execute(
    "detect moustaches within images of faces with a high degree of accuracy.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "detect moustaches within images of faces.", 2)

# Remove the machine-learning sentence that followed "... MRI brain scans."
$d.Content.Find.Execute(
    ". In a real implementation, machine learning would most likely be used in conjunction with the mathematical technique, however since we are mainly focusing on the mathematics, our detector did not rely on any machine learning techniques. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ". ", 2)
